# Correct the mis-keyed product code in C4 (missing leading "1":
# 1321607 -> 11321607) and leave the selection on the cell the user
# was last working with (C17), matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 11321607

$ws.Range("C17").Select() | Out-Null
